# feat: added generate strategy logic
#
# Collapses the various "{{placeholder}}" text fields that were split across
# multiple <a:r> runs back into a single run each, and nudges one shape's
# horizontal offset.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-ShapeText {
    param(
        [__ComObject]$Shape,
        [string]$NewText
    )
    # Assigning the exact same concatenated text the runs already display is
    # a no-op for the run-layout (it stays split across runs), so first push
    # a throwaway value to collapse everything into a single run, then set
    # the real text. The surviving run keeps the formatting of the original
    # first run, matching how PowerPoint merges runs on text assignment.
    $Shape.TextFrame.TextRange.Text = "__tmp__"
    $Shape.TextFrame.TextRange.Text = $NewText
}

# {{icon1}} - was split into "{{" / "icon" / "1}}"
Set-ShapeText -Shape $s.Shapes.Item(3) -NewText "{{icon1}}"

# {{overviewt}} -> {{overview}}
$s.Shapes.Item(11).TextFrame.TextRange.Text = "{{overview}}"

# {{icon2}} - was split into "{{" / "icon" / "2}}"
Set-ShapeText -Shape $s.Shapes.Item(22) -NewText "{{icon2}}"

# {{icon3}} - was split into "{{" / "icon" / "3}}"
Set-ShapeText -Shape $s.Shapes.Item(23) -NewText "{{icon3}}"

# {{icon5}} - was split into "{{" / "icon" / "5}}"
Set-ShapeText -Shape $s.Shapes.Item(24) -NewText "{{icon5}}"

# {{icon4}} - was split into "{{" / "icon" / "4}}"
Set-ShapeText -Shape $s.Shapes.Item(25) -NewText "{{icon4}}"

# {{cash}} - was split into "{{c" / "ash}}"
Set-ShapeText -Shape $s.Shapes.Item(42) -NewText "{{cash}}"

# {{comparativeIndex}} - was split into "{{c" / "omparativeIndex}}"
Set-ShapeText -Shape $s.Shapes.Item(45) -NewText "{{comparativeIndex}}"

# {{sector}} - was split into "{{s" / "ector}}"
Set-ShapeText -Shape $s.Shapes.Item(46) -NewText "{{sector}}"

# {{fundStructure}} - was split into "{{f" / "undStructure}}"
Set-ShapeText -Shape $s.Shapes.Item(47) -NewText "{{fundStructure}}"

# {{ocf}} - was split into "{{ocf" / "}}"
Set-ShapeText -Shape $s.Shapes.Item(49) -NewText "{{ocf}}"

# Shift the portfolioManager2 textbox right (x: 1592500 -> 1643275 EMU)
$shp = $s.Shapes.Item(53)
$shp.Left = 1643275 / 12700
